$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.4401443333333333
$ws.Range("H2").Value = 1.320433
$ws.Range("I2").Value = 0.08872023983698565
$ws.Range("J2").Value = 0.08872023983698567
$ws.Range("M2").Value = 0.4401443333333333
$ws.Range("N2").Value = 1.320433
$ws.Range("O2").Value = 0.08872023983698565
$ws.Range("P2").Value = 0.08872023983698567
$ws.Range("Q2").Value = 0.1937270341654444
$ws.Range("R2").Value = 1.743543307489
$ws.Range("S2").Value = 0.007871280956732256
$ws.Range("T2").Value = 0.007871280956732259
$ws.Range("G3").Value = 0.4401443333333333
$ws.Range("H3").Value = 1.320433
$ws.Range("I3").Value = 0.08872023983698565
$ws.Range("J3").Value = 0.08872023983698567
$ws.Range("O3").Value = 0.6680231260820571
$ws.Range("P3").Value = 0.6680231260820572
$ws.Range("Q3").Value = 1.458676613223667
$ws.Range("R3").Value = 13.128089519013
$ws.Range("S3").Value = 0.05926717196265301
$ws.Range("T3").Value = 0.05926717196265303
$ws.Range("G4").Value = 0.4401443333333333
$ws.Range("H4").Value = 1.320433
$ws.Range("I4").Value = 0.08872023983698565
$ws.Range("J4").Value = 0.08872023983698567
$ws.Range("M4").Value = 1.206805
$ws.Range("N4").Value = 3.620415
$ws.Range("O4").Value = 0.2432566340809571
$ws.Range("P4").Value = 0.2432566340809572
$ws.Range("Q4").Value = 0.5311683821883334
$ws.Range("R4").Value = 4.780515439695001
$ws.Range("S4").Value = 0.02158178691760038
$ws.Range("T4").Value = 0.02158178691760038
$ws.Range("I5").Value = 0.6680231260820571
$ws.Range("J5").Value = 0.6680231260820572
$ws.Range("M5").Value = 0.4401443333333333
$ws.Range("N5").Value = 1.320433
$ws.Range("O5").Value = 0.08872023983698565
$ws.Range("P5").Value = 0.08872023983698567
$ws.Range("Q5").Value = 1.458676613223667
$ws.Range("R5").Value = 13.128089519013
$ws.Range("S5").Value = 0.05926717196265301
$ws.Range("T5").Value = 0.05926717196265303
$ws.Range("I6").Value = 0.6680231260820571
$ws.Range("J6").Value = 0.6680231260820572
$ws.Range("O6").Value = 0.6680231260820571
$ws.Range("P6").Value = 0.6680231260820572
$ws.Range("S6").Value = 0.4462548969804439
$ws.Range("T6").Value = 0.446254896980444
$ws.Range("I7").Value = 0.6680231260820571
$ws.Range("J7").Value = 0.6680231260820572
$ws.Range("M7").Value = 1.206805
$ws.Range("N7").Value = 3.620415
$ws.Range("O7").Value = 0.2432566340809571
$ws.Range("P7").Value = 0.2432566340809572
$ws.Range("Q7").Value = 3.999456762035001
$ws.Range("R7").Value = 35.995110858315
$ws.Range("S7").Value = 0.16250105713896
$ws.Range("T7").Value = 0.1625010571389601
$ws.Range("G8").Value = 1.206805
$ws.Range("H8").Value = 3.620415
$ws.Range("I8").Value = 0.2432566340809571
$ws.Range("J8").Value = 0.2432566340809572
$ws.Range("M8").Value = 0.4401443333333333
$ws.Range("N8").Value = 1.320433
$ws.Range("O8").Value = 0.08872023983698565
$ws.Range("P8").Value = 0.08872023983698567
$ws.Range("Q8").Value = 0.5311683821883334
$ws.Range("R8").Value = 4.780515439695001
$ws.Range("S8").Value = 0.02158178691760038
$ws.Range("T8").Value = 0.02158178691760038
$ws.Range("G9").Value = 1.206805
$ws.Range("H9").Value = 3.620415
$ws.Range("I9").Value = 0.2432566340809571
$ws.Range("J9").Value = 0.2432566340809572
$ws.Range("O9").Value = 0.6680231260820571
$ws.Range("P9").Value = 0.6680231260820572
$ws.Range("Q9").Value = 3.999456762035001
$ws.Range("R9").Value = 35.995110858315
$ws.Range("S9").Value = 0.16250105713896
$ws.Range("T9").Value = 0.1625010571389601
$ws.Range("G10").Value = 1.206805
$ws.Range("H10").Value = 3.620415
$ws.Range("I10").Value = 0.2432566340809571
$ws.Range("J10").Value = 0.2432566340809572
$ws.Range("M10").Value = 1.206805
$ws.Range("N10").Value = 3.620415
$ws.Range("O10").Value = 0.2432566340809571
$ws.Range("P10").Value = 0.2432566340809572
$ws.Range("Q10").Value = 1.456378308025
$ws.Range("R10").Value = 13.107404772225
$ws.Range("S10").Value = 0.05917379002439668
$ws.Range("T10").Value = 0.05917379002439669
